# Fruta / hortaliza, semanal
# New weekly price record for Arandano (blue) - Vega Monumental Concepcion.
# Insert a brand-new data row right above the current row 26, pushing the
# existing rows 26-50 down to 27-51, then populate the new row with its data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 26 (shifts old rows 26:50 down to 27:51)
$ws.Rows.Item(26).Insert()

# Populate the new row 26 with the new weekly record
$ws.Range("A26").Value = 11
$ws.Range("B26").Value = "Vega Monumental Concepción"
$ws.Range("C26").Value = "Bíobío"
$ws.Range("D26").Value = 44159
$ws.Range("E26").Value = 8
$ws.Range("F26").Value = "Fruta"
$ws.Range("G26").Value = 100101
$ws.Range("H26").Value = "Berries"
$ws.Range("I26").Value = 100101001
$ws.Range("J26").Value = "Arándano (blue)"
$ws.Range("K26").Value = "Sin especificar"
$ws.Range("L26").Value = "Primera"
$ws.Range("M26").Value = 100
$ws.Range("N26").Value = 6000
$ws.Range("O26").Value = 6000
$ws.Range("P26").Value = 6000
$ws.Range("Q26").Value = "$/bandeja 2 kilos"
$ws.Range("R26").Value = "Provincia de Curicó"
$ws.Range("S26").Value = 3000
$ws.Range("T26").Value = 2
